# Jononi Scripts.xlsx - numeric and radio button add
# 1. Rename sheet "MPerformanceReviewMeeting" -> "M_Performance_Review_Meeting"
# 2. On "S Awareness Session", change D35 datatype value from "PDF" to "pdf"

$wb = $excel.ActiveWorkbook

# Rename the performance-review sheet to use underscores.
$perfSheet = $wb.Worksheets.Item("MPerformanceReviewMeeting")
$perfSheet.Name = "M_Performance_Review_Meeting"

# Correct the datatype value for the PDF upload question (row 35) on the
# Awareness Session sheet: "PDF" -> "pdf".
$awarenessSheet = $wb.Worksheets.Item("S Awareness Session")
$awarenessSheet.Range("D35").Value = "pdf"

# Leave this as the active sheet/cell, matching the edited workbook's view.
$awarenessSheet.Activate() | Out-Null
$awarenessSheet.Range("D35").Select() | Out-Null
